$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6971663236618042
$ws.Range("B1").Value = 3.750109434127808
$ws.Range("C1").Value = 5.583365917205811
$ws.Range("D1").Value = 1.248443961143494
$ws.Range("E1").Value = 0.71803218126297
